$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B2" = 1.129710549422754
    "D2" = 0.1439484244774292
    "E2" = 0.1709925828979948
    "F2" = 1.913762155300162
    "G2" = 1.424422803818146
    "H2" = 1.253037334986175
    "I2" = 1.026254946460977
    "J2" = 0.2333503628092757
    "L2" = 0.5719242724773324
    "M2" = 0.397124536698243
    "N2" = 1.705392521894602
    "B3" = 1.074200315205076
    "D3" = 0.138153360017057
    "E3" = 0.1630573364092456
    "F3" = 1.902138619123846
    "G3" = 1.394602647724838
    "H3" = 1.245767170004882
    "I3" = 1.040115062785878
    "J3" = 0.2216599351786499
    "L3" = 0.5369131452780209
    "M3" = 0.3754618069485787
    "N3" = 1.71702300785995
    "B4" = 1.040463292841736
    "D4" = 0.1345519039550496
    "E4" = 0.1581391367878098
    "F4" = 1.896237796578859
    "G4" = 1.377362996401672
    "H4" = 1.242062647602154
    "I4" = 1.04915934569064
    "J4" = 0.2144233159528284
    "L4" = 0.5156206934085503
    "M4" = 0.3622924937205525
    "N4" = 1.724786803803724
    "B5" = 1.026803131409281
    "D5" = 0.1330732933918881
    "E5" = 0.1561233434492664
    "F5" = 1.894143248181621
    "G5" = 1.370605587835797
    "H5" = 1.24074359113385
    "I5" = 1.052979143706674
    "J5" = 0.2114595995698068
    "L5" = 0.5069953875702993
    "M5" = 0.3569592133494268
    "N5" = 1.728107436477522
    "B6" = 1.024540208268519
    "D6" = 0.1328271044141971
    "E6" = 0.1557879215687414
    "F6" = 1.893814156993983
    "G6" = 1.369499670762906
    "H6" = 1.240536061329578
    "I6" = 1.05362151862067
    "J6" = 0.2109665871540329
    "L6" = 0.5055662785304094
    "M6" = 0.3560756449197271
    "N6" = 1.728668304403378
    "B7" = 1.040278709935876
    "D7" = 0.1345320075264524
    "E7" = 0.1581119980575139
    "F7" = 1.896208294229424
    "G7" = 1.377270780519325
    "H7" = 1.242044087271807
    "I7" = 1.049210317710255
    "J7" = 0.2143834058797438
    "L7" = 0.5155041606778923
    "M7" = 0.3622204320435358
    "N7" = 1.724830951707979
    "B8" = 1.1104991800525
    "D8" = 0.1419591982531898
    "E8" = 0.1682660163710139
    "F8" = 1.909497358052349
    "G8" = 1.413917999278709
    "H8" = 1.250372730427102
    "I8" = 1.030922982750308
    "J8" = 0.2293316800398202
    "L8" = 0.5598100415368492
    "M8" = 0.3896280112354091
    "N8" = 1.709273638701056
    "B9" = 1.250922290934
    "D9" = 0.1561861214397595
    "E9" = 0.1878162413001832
    "F9" = 1.945402487263976
    "G9" = 1.494332703216458
    "H9" = 1.272750868124007
    "I9" = 0.9993049111922403
    "J9" = 0.2581806622534515
    "L9" = 0.6483167224948829
    "M9" = 0.4444131167562873
    "N9" = 1.683695194549585
    "B10" = 1.355723670794134
    "D10" = 0.1664411654595881
    "E10" = 0.2019641710143461
    "F10" = 1.977842360598999
    "G10" = 1.558716158103863
    "H10" = 1.292909830617958
    "I10" = 0.9786688003388342
    "J10" = 0.2790959494226541
    "L10" = 0.7143398781350925
    "M10" = 0.485293671741502
    "N10" = 1.667893756055307
    "B11" = 1.403750277914014
    "D11" = 0.1710654465030927
    "E11" = 0.2083547580269638
    "F11" = 1.993929175902537
    "G11" = 1.589177517364789
    "H11" = 1.302895213528387
    "I11" = 0.9698455784568054
    "J11" = 0.288550852160796
    "L11" = 0.744594323536063
    "M11" = 0.5040276895206546
    "N11" = 1.661351950973199
    "B12" = 1.42198667030209
    "D12" = 0.1728108100522405
    "E12" = 0.2107682468441965
    "F12" = 2.00021301452287
    "G12" = 1.600882601025233
    "H12" = 1.306794137091742
    "I12" = 0.966585740960614
    "J12" = 0.2921226345225847
    "L12" = 0.7560825863351965
    "M12" = 0.5111413811053893
    "N12" = 1.658967472021601
    "B13" = 1.418056937321467
    "D13" = 0.1724351687447836
    "E13" = 0.2102487467849556
    "F13" = 1.998851118917429
    "G13" = 1.598354118129208
    "H13" = 1.305949192692594
    "I13" = 0.9672841853349361
    "J13" = 0.2913537693459034
    "L13" = 0.7536069812474011
    "M13" = 0.5096084552522413
    "N13" = 1.659476889509676
    "B14" = 1.40524960449244
    "D14" = 0.1712091532193654
    "E14" = 0.2085534470114254
    "F14" = 1.99444229574533
    "G14" = 1.590137086241299
    "H14" = 1.303213619121664
    "I14" = 0.9695757587096168
    "J14" = 0.2888448769852658
    "L14" = 0.7455388378584189
    "M14" = 0.504612547080562
    "N14" = 1.661153920416581
    "B15" = 1.397411198294719
    "D15" = 0.1704574377269523
    "E15" = 0.2075141831752561
    "F15" = 1.991766809448805
    "G15" = 1.58512609994068
    "H15" = 1.301553340666601
    "I15" = 0.9709900100843036
    "J15" = 0.2873069895394309
    "L15" = 0.7406009743481832
    "M15" = 0.5015549456911401
    "N15" = 1.662193225133194
    "B16" = 1.352592028126196
    "D16" = 0.1661381467591667
    "E16" = 0.201545621192281
    "F16" = 1.976817883603758
    "G16" = 1.556749155711543
    "H16" = 1.292273702857955
    "I16" = 0.9792567886909076
    "J16" = 0.2784768490511738
    "L16" = 0.7123671008197334
    "M16" = 0.4840721033416671
    "N16" = 1.668334268724649
    "B17" = 1.325186451103775
    "D17" = 0.1634780352940055
    "E17" = 0.1978725058419499
    "F17" = 1.96798841636847
    "G17" = 1.539642211503292
    "H17" = 1.286790021145634
    "I17" = 0.9844728394999329
    "J17" = 0.2730445739052669
    "L17" = 0.6951028538268815
    "M17" = 0.4733819476396448
    "N17" = 1.67226701238009
    "B18" = 1.309456664013794
    "D18" = 0.1619441632985001
    "E18" = 0.1957555536278832
    "F18" = 1.963035049069489
    "G18" = 1.529913076688189
    "H18" = 1.283712631920167
    "I18" = 0.9875260724179249
    "J18" = 0.2699144672820353
    "L18" = 0.6851936456793055
    "M18" = 0.467246183670845
    "N18" = 1.674589870386271
    "B19" = 1.30413655237453
    "D19" = 0.1614241565408747
    "E19" = 0.1950380548982409
    "F19" = 1.961379384478448
    "G19" = 1.526637866590306
    "H19" = 1.282683835929333
    "I19" = 0.9885689586633717
    "J19" = 0.2688537047772996
    "L19" = 0.6818421216497086
    "M19" = 0.4651709467677065
    "N19" = 1.675386806586388
    "B20" = 1.328100394162334
    "D20" = 0.1637616059428808
    "E20" = 0.1982639572509015
    "F20" = 1.968915373544988
    "G20" = 1.54145184675761
    "H20" = 1.287365828937226
    "I20" = 0.9839120851675958
    "J20" = 0.2736234291630524
    "L20" = 0.6969385176561502
    "M20" = 0.474518596207794
    "N20" = 1.671842068953069
    "B21" = 1.409010084894533
    "D21" = 0.1715694188721528
    "E21" = 0.2090515734233662
    "F21" = 1.995732052975015
    "G21" = 1.59254600234533
    "H21" = 1.304013925810182
    "I21" = 0.9689004594200519
    "J21" = 0.2895820325990854
    "L21" = 0.7479077913643835
    "M21" = 0.5060794381550977
    "N21" = 1.660658819712808
    "B22" = 1.462178887441723
    "D22" = 0.1766388066916704
    "E22" = 0.2160641064893127
    "F22" = 2.014378494211527
    "G22" = 1.626930676222599
    "H22" = 1.315580553673442
    "I22" = 0.9595636075458245
    "J22" = 0.2999618600839682
    "L22" = 0.7814030770724969
    "M22" = 0.5268199877239539
    "N22" = 1.65389052031999
    "B23" = 1.433775473173569
    "D23" = 0.1739362034135041
    "E23" = 0.212324831487031
    "F23" = 2.004323741985104
    "G23" = 1.608487739128492
    "H23" = 1.309344283559767
    "I23" = 0.9645034174163634
    "J23" = 0.2944265354725673
    "L23" = 0.763509215811041
    "M23" = 0.5157400298963992
    "N23" = 1.657453482699367
    "B24" = 1.326782919313018
    "D24" = 0.1636334178007957
    "E24" = 0.1980869983622995
    "F24" = 1.968495913726656
    "G24" = 1.540633380973702
    "H24" = 1.287105271862089
    "I24" = 0.984165432589748
    "J24" = 0.2733617505661528
    "L24" = 0.6961085634541178
    "M24" = 0.4740046857607041
    "N24" = 1.672033993138072
    "B25" = 1.212645561671934
    "D25" = 0.1523726541197732
    "E25" = 0.1825656119314232
    "F25" = 1.934629189134085
    "G25" = 1.471653767640902
    "H25" = 1.266046459764965
    "I25" = 1.007403460107454
    "J25" = 0.2504256605021453
    "L25" = 0.6241987755155094
    "M25" = 0.4294814347178786
    "N25" = 1.690088598177645
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

Write-Output "Updated $($updates.Count) cells"